# New changes, included calculate button
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "Delivery Time"
$ws.Range("G1").Value = "Order Details"

# --- Row 2 ---
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "12"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "12.00"
$ws.Range("D2").Value = "Delivery"
$ws.Range("E2").Value = "Not Paid"
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "Rice: 1 | `nBhaji:   | `nBhakari:   | `nVaran: 1 | `n"

# --- Row 3 ---
$ws.Range("A3").Value = "b"
$ws.Range("B3").ClearContents()
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "15.00"
$ws.Range("D3").Value = "Take Away"
$ws.Range("E3").Value = "Not Paid"
$ws.Range("F3").Value = "Lomesh"
$ws.Range("G3").Value = "Rice: 0 | `nBhaji:   | `nBhakari:   | `nVaran: 0 | `n"

# --- Row 4 ---
$ws.Range("A4").Value = "d"
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = 60
$ws.Range("D4").Value = "Delivery"
$ws.Range("E4").Value = "Not Paid"
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "Rice: 1 | `nBhaji:   | `nBhakari:   | `nVaran: 1 | `n"

# --- Row 5 ---
$ws.Range("A5").Value = "f"
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = 110
$ws.Range("D5").Value = "Take Away"
$ws.Range("E5").Value = "Not Paid"
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = "Rice: 1 | `nBhaji: 1  | `nBhakari: 3  | `nVaran: 2 | `n"

# --- Row 6 (new) ---
$ws.Range("A6").Value = "f"
$ws.Range("C6").Value = 80
$ws.Range("D6").Value = "Take Away"
$ws.Range("E6").Value = "Not Paid"
$ws.Range("G6").Value = "Rice: 0 | `nBhaji: 1  | `nBhakari: 3  | `nVaran: 1 | `n"

# --- Row 7 (new) ---
$ws.Range("A7").Value = "f"
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = "Take Away"
$ws.Range("E7").Value = "Not Paid"
$ws.Range("G7").Value = "Rice: 0 | `nBhaji: 1  | `nBhakari: 3  | `nVaran: 0 | `n"

# --- Row 8 (new) ---
$ws.Range("A8").Value = "f"
$ws.Range("C8").Value = 60
$ws.Range("D8").Value = "Take Away"
$ws.Range("E8").Value = "Not Paid"
$ws.Range("G8").Value = "Rice: 0 | `nBhaji: 1special  | `nBhakari: 3matar paneer  | `nVaran: 0 | `n"

# --- Row 9 (new) ---
$ws.Range("A9").Value = "g"
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = "Take Away"
$ws.Range("E9").Value = "Not Paid"
$ws.Range("F9").Value = "Rajas"
$ws.Range("G9").Value = "Rice: 1 | `nBhaji: 1 matar  | `nBhakari: 1 bhendi  | `nVaran: 1 | `n"

# --- Row 10 (new) ---
$ws.Range("A10").Value = "aa"
$ws.Range("C10").Value = 60
$ws.Range("D10").Value = "Delivery"
$ws.Range("E10").Value = "Not Paid"
$ws.Range("G10").Value = "Rice: 1 | `nBhaji: 1   | `nBhakari: 1   | `nVaran: 1 | `n"

# --- Column widths (autofit to match bestFit columns) ---
$ws.Columns("A:G").AutoFit()

# --- Selection / view state ---
$ws.Range("C6").Select()

# --- Recalculation button / settings (calcId reset, matches "included calculate button") ---
$excel.CalculateFullRebuild()
